# dropped _ from entries
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A35").Value = "sabertoothcat"
$ws.Range("A43").Value = "trex"
$ws.Range("A59").Value = "ringtailed"
